$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to a handful of already-reported days (column C = "Nb
# nouveaux cas positifs"); the cumulative column B, and the downstream
# H/J/K helper formulas, ripple automatically on recalculation. ---
$ws.Range("C321").Value = 197
$ws.Range("C335").Value = 154
$ws.Range("C370").Value = 89
$ws.Range("C371").Value = 69
$ws.Range("C372").Value = 59
$ws.Range("C373").Value = 38

# L372, L373 and M373 were stored as text ("1"/"0") because the column is
# text-formatted ("@"); flip to General just long enough to store a real
# number, then restore the original (text) display format so the cell
# style index is unchanged.
$fmt = $ws.Range("L372").NumberFormat
$ws.Range("L372").NumberFormat = "General"
$ws.Range("L372").Value = 1
$ws.Range("L372").NumberFormat = $fmt

$fmt = $ws.Range("L373").NumberFormat
$ws.Range("L373").NumberFormat = "General"
$ws.Range("L373").Value = 0
$ws.Range("L373").NumberFormat = $fmt

$fmt = $ws.Range("M373").NumberFormat
$ws.Range("M373").NumberFormat = "General"
$ws.Range("M373").Value = 0
$ws.Range("M373").NumberFormat = $fmt

# --- New day entered: row 374 (05.03.2021) ---
$ws.Range("C374").Value = 8
$ws.Range("E374").Value = 9
$ws.Range("F374").Value = 7
$ws.Range("G374").Value = 22

$fmt = $ws.Range("L374").NumberFormat
$ws.Range("L374").NumberFormat = "General"
$ws.Range("L374").Value = 0
$ws.Range("L374").NumberFormat = $fmt

$fmt = $ws.Range("M374").NumberFormat
$ws.Range("M374").NumberFormat = "General"
$ws.Range("M374").Value = 0
$ws.Range("M374").NumberFormat = $fmt

# --- View state: scroll the frozen pane back to the top and move the
# selection back to A2 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$null = $ws.Range("A2").Select()
